$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)

# --- 1. Resize / reposition the content placeholder -----------------------
# off  x=618067  y=1027906 -> y=1027907
# ext cx=10515600 cy=5726641 -> cy=5830094
# (Shape.Top/Left/Width/Height are expressed in points; 914400 EMU = 72 pt)
$sh.Top    = 1027907 / 914400 * 72
$sh.Height = 5830094 / 914400 * 72

$tr = $sh.TextFrame.TextRange

# --- 2. Merge the RF07 paragraph's two runs into a single run -------------
$rf07 = $tr.Paragraphs(8, 1)
$rf07Text = "RF07 – O sistema deve criar um dashboard mostrando os dados para o fabricante com prazo de entrega, e dados gerais sobre os pedidos que foram feitos e o tempo de envio para o correio."
# Collapse to a placeholder first so the merge isn't re-split back into the
# original two run "seats" when the final text is written.
$rf07.Text = "X"
$rf07 = $tr.Paragraphs(8, 1)
$rf07.Text = $rf07Text

# --- 3. Give the "Não-funcionais:" paragraph 120% line spacing ------------
$naoFuncionais = $tr.Paragraphs(10, 1)
$naoFuncionais.ParagraphFormat.SpaceWithin = 1.2
